# envSampCleanMapping.xlsx update
# Refreshes the sample mapping table (rows that previously held placeholder
# "Sample <id>" / NA / None values for the lab-derived SealCoat samples) with
# the corrected sample_matrix / technology / client / location data, and
# turns on AutoFilter for the data range now that the table is finalized.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Gulf Seal No. 253 (Henry Company, CT-2) composite/synthetic rows ---
$ws.Range("C21").Value = "Synthetic"
$ws.Range("D21").Value = "Composite"

$ws.Range("C22").Value = "Synthetic"
$ws.Range("D22").Value = "Composite"

$ws.Range("C23").Value = "Synthetic"
$ws.Range("D23").Value = "Composite"

# --- Henry Seal No. 532 (Henry Company, AS) composite/synthetic rows ---
$ws.Range("C38").Value = "Synthetic"
$ws.Range("D38").Value = "Composite"

$ws.Range("C39").Value = "Synthetic"
$ws.Range("D39").Value = "Composite"

$ws.Range("C40").Value = "Synthetic"
$ws.Range("D40").Value = "Composite"

# --- Row 49: was a lab-only placeholder ("Sample 101"); now resolved to the
#     Tarconite (Neyra Industries) SealCoat_1.6 h B sample ---
$ws.Range("C49").Value = "Synthetic"
$ws.Range("D49").Value = "Composite"
$ws.Range("F49").Value = "Tarconite Neyra Industries Inc."
$ws.Range("G49").Value = "SealCoat_1.6 h B"
$ws.Range("J49").Value = "Tarconite Neyra Industries Inc."
$ws.Range("L49").Value = "CT-1 on pavement"
$ws.Range("M49").Value = "Tarconite Neyra Industries Inc."
$ws.Range("N49").Value = "SealCoat_1.6 h B"
$ws.Range("O49").Value = "Tarconite Neyra Industries Inc."

# --- Row 50: "Sample 103" -> SealCoat_1.6 h C (Tarconite) ---
$ws.Range("C50").Value = "Synthetic"
$ws.Range("D50").Value = "Composite"
$ws.Range("F50").Value = "Tarconite Neyra Industries Inc."
$ws.Range("G50").Value = "SealCoat_1.6 h C"
$ws.Range("J50").Value = "Tarconite Neyra Industries Inc."
$ws.Range("L50").Value = "CT-1 on pavement"
$ws.Range("M50").Value = "Tarconite Neyra Industries Inc."
$ws.Range("N50").Value = "SealCoat_1.6 h C"
$ws.Range("O50").Value = "Tarconite Neyra Industries Inc."

# --- Row 51: "Sample 1035" -> SealCoat_AS A (Henry Seal No. 532) ---
$ws.Range("C51").Value = "Synthetic"
$ws.Range("D51").Value = "Composite"
$ws.Range("F51").Value = "Henry Seal No. 532, Henry Company"
$ws.Range("G51").Value = "SealCoat_AS A"
$ws.Range("J51").Value = "Henry Seal No. 532, Henry Company"
$ws.Range("M51").Value = "Henry Seal No. 532, Henry Company"
$ws.Range("N51").Value = "SealCoat_AS A"
$ws.Range("O51").Value = "Henry Seal No. 532, Henry Company"

# --- Row 52: "Sample 1514" -> SealCoat_CT-1 A (Tarconite) ---
$ws.Range("C52").Value = "Synthetic"
$ws.Range("D52").Value = "Composite"
$ws.Range("F52").Value = "Tarconite Neyra Industries Inc."
$ws.Range("G52").Value = "SealCoat_CT-1 A"
$ws.Range("J52").Value = "Tarconite Neyra Industries Inc."
$ws.Range("L52").Value = "CT-1"
$ws.Range("M52").Value = "Tarconite Neyra Industries Inc."
$ws.Range("N52").Value = "SealCoat_CT-1 A"
$ws.Range("O52").Value = "Tarconite Neyra Industries Inc."

# --- Row 53: "Sample 1518" -> SealCoat_CT-2 A (Gulf Seal No. 253) ---
$ws.Range("C53").Value = "Synthetic"
$ws.Range("D53").Value = "Composite"
$ws.Range("F53").Value = "Gulf Seal No. 253 Henry Company"
$ws.Range("G53").Value = "SealCoat_CT-2 A"
$ws.Range("J53").Value = "Gulf Seal No. 253 Henry Company"
$ws.Range("L53").Value = "CT-2"
$ws.Range("M53").Value = "Gulf Seal No. 253 Henry Company"
$ws.Range("N53").Value = "SealCoat_CT-2 A"
$ws.Range("O53").Value = "Gulf Seal No. 253 Henry Company"

# --- Row 54: matrix/technology now known (stays "Sample 2791" lab sample) ---
$ws.Range("C54").Value = "Synthetic"
$ws.Range("D54").Value = "Composite"

# Table is final now -- turn on AutoFilter over the whole data range.
$ws.Range("A1:O132").AutoFilter()

# AutoFilter should register the (hidden, sheet-scoped) _FilterDatabase name,
# matching what Excel itself writes out when a filter is turned on.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=envSampCleanMapping!`$A`$1:`$O`$132")
$filterName.Visible = $false

# Leave the selection roughly where the authoring user left it.
$ws.Range("C54").Select()
